$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = 'Datos actualizados a 16 de Octubre de 2020 a las 18:23'

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 8230562
$ws.Cells.Item(4,3).Value = 14247
$ws.Cells.Item(4,4).Value = 5329170
$ws.Cells.Item(4,5).Value = 2678411
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 264
$ws.Cells.Item(4,8).Value = 222981

# Row 5: India
$ws.Cells.Item(5,2).Value = 7420529
$ws.Cells.Item(5,3).Value = 55020
$ws.Cells.Item(5,4).Value = 6509841
$ws.Cells.Item(5,5).Value = 797774
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 768
$ws.Cells.Item(5,8).Value = 112914

# Row 6: Brasil
$ws.Cells.Item(6,2).Value = 5176524
$ws.Cells.Item(6,3).Value = 5528
$ws.Cells.Item(6,4).Value = 4599446
$ws.Cells.Item(6,5).Value = 424437
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 128
$ws.Cells.Item(6,8).Value = 152641

# Row 15: Reino Unido
$ws.Cells.Item(15,2).Value = 689257
$ws.Cells.Item(15,3).Value = 15650
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 136
$ws.Cells.Item(15,8).Value = 43429

# Row 19: Italia
$ws.Cells.Item(19,1).Value = 'Italia'
$ws.Cells.Item(19,2).Value = 391611
$ws.Cells.Item(19,3).Value = 10010
$ws.Cells.Item(19,4).Value = 247872
$ws.Cells.Item(19,5).Value = 107312
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = 55
$ws.Cells.Item(19,8).Value = 36427

# Row 20: Banglades
$ws.Cells.Item(20,1).Value = 'Banglades'
$ws.Cells.Item(20,2).Value = 386086
$ws.Cells.Item(20,3).Value = 1527
$ws.Cells.Item(20,4).Value = 300738
$ws.Cells.Item(20,5).Value = 79725
$ws.Cells.Item(20,6).Value = 0
$ws.Cells.Item(20,7).Value = 15
$ws.Cells.Item(20,8).Value = 5623

# Row 21: Alemania
$ws.Cells.Item(21,2).Value = 354643
$ws.Cells.Item(21,3).Value = 5827
$ws.Cells.Item(21,4).Value = 284600
$ws.Cells.Item(21,5).Value = 60212
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 21
$ws.Cells.Item(21,8).Value = 9831

# Row 27: Israel
$ws.Cells.Item(27,2).Value = 301896
$ws.Cells.Item(27,3).Value = 1695
$ws.Cells.Item(27,4).Value = 262503
$ws.Cells.Item(27,5).Value = 37252
$ws.Cells.Item(27,6).Value = 0
$ws.Cells.Item(27,7).Value = 14
$ws.Cells.Item(27,8).Value = 2141

# Row 30: Canada
$ws.Cells.Item(30,2).Value = 193497
$ws.Cells.Item(30,3).Value = 1765
$ws.Cells.Item(30,4).Value = 163248
$ws.Cells.Item(30,5).Value = 20528
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 22
$ws.Cells.Item(30,8).Value = 9721

# Row 35: Chequia
$ws.Cells.Item(35,1).Value = 'Chequia'
$ws.Cells.Item(35,2).Value = 154675
$ws.Cells.Item(35,3).Value = 5665
$ws.Cells.Item(35,4).Value = 66005
$ws.Cells.Item(35,5).Value = 87398
$ws.Cells.Item(35,6).Value = 0
$ws.Cells.Item(35,7).Value = 42
$ws.Cells.Item(35,8).Value = 1272

# Row 36: Ecuador
$ws.Cells.Item(36,1).Value = 'Ecuador'
$ws.Cells.Item(36,2).Value = 150360
$ws.Cells.Item(36,3).Value = 0
$ws.Cells.Item(36,4).Value = 128134
$ws.Cells.Item(36,5).Value = 9920
$ws.Cells.Item(36,6).Value = 0
$ws.Cells.Item(36,7).Value = 0
$ws.Cells.Item(36,8).Value = 12306

# Row 38: Catar
$ws.Cells.Item(38,2).Value = 128992
$ws.Cells.Item(38,3).Value = 189
$ws.Cells.Item(38,4).Value = 126006
$ws.Cells.Item(38,5).Value = 2764
$ws.Cells.Item(38,6).Value = 0
$ws.Cells.Item(38,7).Value = 0
$ws.Cells.Item(38,8).Value = 222

# Row 41: Republica Dominicana
$ws.Cells.Item(41,2).Value = 120450
$ws.Cells.Item(41,3).Value = 384
$ws.Cells.Item(41,4).Value = 96883
$ws.Cells.Item(41,5).Value = 21375
$ws.Cells.Item(41,6).Value = 0
$ws.Cells.Item(41,7).Value = 3
$ws.Cells.Item(41,8).Value = 2192

# Row 42: Kuwait
$ws.Cells.Item(42,2).Value = 114744
$ws.Cells.Item(42,3).Value = 729
$ws.Cells.Item(42,4).Value = 106495
$ws.Cells.Item(42,5).Value = 7559
$ws.Cells.Item(42,6).Value = 0
$ws.Cells.Item(42,7).Value = 6
$ws.Cells.Item(42,8).Value = 690

# Row 79: Jordania
$ws.Cells.Item(79,1).Value = 'Jordania'
$ws.Cells.Item(79,2).Value = 34548
$ws.Cells.Item(79,3).Value = 1539
$ws.Cells.Item(79,4).Value = 6692
$ws.Cells.Item(79,5).Value = 27546
$ws.Cells.Item(79,6).Value = 0
$ws.Cells.Item(79,7).Value = 28
$ws.Cells.Item(79,8).Value = 310

# Row 80: Dinamarca
$ws.Cells.Item(80,1).Value = 'Dinamarca'
$ws.Cells.Item(80,2).Value = 34441
$ws.Cells.Item(80,3).Value = 418
$ws.Cells.Item(80,4).Value = 28551
$ws.Cells.Item(80,5).Value = 5213
$ws.Cells.Item(80,6).Value = 0
$ws.Cells.Item(80,7).Value = 0
$ws.Cells.Item(80,8).Value = 677

# Row 81: Birmania
$ws.Cells.Item(81,1).Value = 'Birmania'
$ws.Cells.Item(81,2).Value = 33488
$ws.Cells.Item(81,3).Value = 1137
$ws.Cells.Item(81,4).Value = 15477
$ws.Cells.Item(81,5).Value = 17212
$ws.Cells.Item(81,6).Value = 0
$ws.Cells.Item(81,7).Value = 34
$ws.Cells.Item(81,8).Value = 799

# Row 88: Grecia
$ws.Cells.Item(88,2).Value = 24450
$ws.Cells.Item(88,3).Value = 503
$ws.Cells.Item(88,4).Value = 9989
$ws.Cells.Item(88,5).Value = 13971
$ws.Cells.Item(88,6).Value = 0
$ws.Cells.Item(88,7).Value = 8
$ws.Cells.Item(88,8).Value = 490

# Row 100: Montenegro
$ws.Cells.Item(100,2).Value = 15281
$ws.Cells.Item(100,3).Value = 273
$ws.Cells.Item(100,4).Value = 10569
$ws.Cells.Item(100,5).Value = 4484
$ws.Cells.Item(100,6).Value = 0
$ws.Cells.Item(100,7).Value = 7
$ws.Cells.Item(100,8).Value = 228

# Row 109: Luxemburgo
$ws.Cells.Item(109,1).Value = 'Luxemburgo'
$ws.Cells.Item(109,2).Value = 10471
$ws.Cells.Item(109,3).Value = 227
$ws.Cells.Item(109,4).Value = 8468
$ws.Cells.Item(109,5).Value = 1870
$ws.Cells.Item(109,6).Value = 0
$ws.Cells.Item(109,7).Value = 0
$ws.Cells.Item(109,8).Value = 133

# Row 110: Tayikistan
$ws.Cells.Item(110,1).Value = 'Tayikistan'
$ws.Cells.Item(110,2).Value = 10414
$ws.Cells.Item(110,3).Value = 40
$ws.Cells.Item(110,4).Value = 9393
$ws.Cells.Item(110,5).Value = 941
$ws.Cells.Item(110,6).Value = 0
$ws.Cells.Item(110,7).Value = 0
$ws.Cells.Item(110,8).Value = 80

# Row 111: Uganda
$ws.Cells.Item(111,1).Value = 'Uganda'
$ws.Cells.Item(111,2).Value = 10334
$ws.Cells.Item(111,3).Value = 217
$ws.Cells.Item(111,4).Value = 6901
$ws.Cells.Item(111,5).Value = 3337
$ws.Cells.Item(111,6).Value = 0
$ws.Cells.Item(111,7).Value = 0
$ws.Cells.Item(111,8).Value = 96

# Row 112: Guayana Francesa
$ws.Cells.Item(112,2).Value = 10239
$ws.Cells.Item(112,3).Value = 6
$ws.Cells.Item(112,4).Value = 9955
$ws.Cells.Item(112,5).Value = 215
$ws.Cells.Item(112,6).Value = 0
$ws.Cells.Item(112,7).Value = 0
$ws.Cells.Item(112,8).Value = 69

# Row 137: Reunion
$ws.Cells.Item(137,2).Value = 4776
$ws.Cells.Item(137,3).Value = 98
$ws.Cells.Item(137,4).Value = 4445
$ws.Cells.Item(137,5).Value = 314
$ws.Cells.Item(137,6).Value = 0
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 17

# Row 149: Principado de Andorra
$ws.Cells.Item(149,1).Value = 'Principado de Andorra'
$ws.Cells.Item(149,2).Value = 3377
$ws.Cells.Item(149,3).Value = 187
$ws.Cells.Item(149,4).Value = 2057
$ws.Cells.Item(149,5).Value = 1261
$ws.Cells.Item(149,6).Value = 0
$ws.Cells.Item(149,7).Value = 0
$ws.Cells.Item(149,8).Value = 59

# Row 150: Mali
$ws.Cells.Item(150,1).Value = 'Mali'
$ws.Cells.Item(150,2).Value = 3368
$ws.Cells.Item(150,3).Value = 0
$ws.Cells.Item(150,4).Value = 2559
$ws.Cells.Item(150,5).Value = 677
$ws.Cells.Item(150,6).Value = 0
$ws.Cells.Item(150,7).Value = 0
$ws.Cells.Item(150,8).Value = 132

# Row 151: Letonia
$ws.Cells.Item(151,1).Value = 'Letonia'
$ws.Cells.Item(151,2).Value = 3204
$ws.Cells.Item(151,3).Value = 148
$ws.Cells.Item(151,4).Value = 1329
$ws.Cells.Item(151,5).Value = 1833
$ws.Cells.Item(151,6).Value = 0
$ws.Cells.Item(151,7).Value = 1
$ws.Cells.Item(151,8).Value = 42

# Row 152: Sudan del Sur
$ws.Cells.Item(152,2).Value = 2817
$ws.Cells.Item(152,3).Value = 10
$ws.Cells.Item(152,4).Value = 1290
$ws.Cells.Item(152,5).Value = 1472
$ws.Cells.Item(152,6).Value = 0
$ws.Cells.Item(152,7).Value = 0
$ws.Cells.Item(152,8).Value = 55
